$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RF)
$ws.Range("B3").Value = 0.901
$ws.Range("C3").Value = 0.893
$ws.Range("D3").Value = 0.108
$ws.Range("E3").Value = 0.329
$ws.Range("F3").Value = 0.24
$ws.Range("G3").Value = 0.97

# Row 4 (NN)
$ws.Range("B4").Value = 0.712
$ws.Range("C4").Value = 0.6879999999999999
$ws.Range("D4").Value = 0.315
$ws.Range("E4").Value = 0.5610000000000001
$ws.Range("F4").Value = 0.435
$ws.Range("G4").Value = 0.907

# Row 5 (RNN)
$ws.Range("B5").Value = 0.636
$ws.Range("C5").Value = 0.62
$ws.Range("D5").Value = 0.397
$ws.Range("E5").Value = 0.63
$ws.Range("F5").Value = 0.474
$ws.Range("G5").Value = 0.832
